$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$logRows = @(
    @("2025-12-25 00:20:51", "Admin", "Login", "login_success", "Role: admin"),
    @("2025-12-25 00:20:52", "Admin", "dashboard", "access_granted", "Opened dashboard page"),
    @("2025-12-25 00:20:53", "Admin", "dashboard", "access_granted", "Opened dashboard page"),
    @("2025-12-25 00:20:54", "Admin", "quotation", "access_granted", "Opened quotation page"),
    @("2025-12-25 00:21:01", "Admin", "invoice", "access_granted", "Opened invoice page"),
    @("2025-12-25 00:21:06", "Admin", "receipt", "access_granted", "Opened receipt page"),
    @("2025-12-25 00:21:08", "Admin", "quotation", "access_granted", "Opened quotation page"),
    @("2025-12-25 00:21:14", "Admin", "invoice", "access_granted", "Opened invoice page"),
    @("2025-12-25 00:21:17", "Admin", "invoice", "access_granted", "Opened invoice page"),
    @("2025-12-25 00:21:17", "Admin", "invoice", "access_granted", "Opened invoice page"),
    @("2025-12-25 00:21:31", "Admin", "invoice", "access_granted", "Opened invoice page"),
    @("2025-12-25 00:21:32", "Admin", "invoice", "access_granted", "Opened invoice page"),
    @("2025-12-25 00:21:33", "Admin", "invoice", "access_granted", "Opened invoice page"),
    @("2025-12-25 00:22:48", "Admin", "quotation", "access_granted", "Opened quotation page"),
    @("2025-12-25 00:35:22", "Admin", "Login", "login_success", "Role: admin"),
    @("2025-12-25 00:35:22", "Admin", "dashboard", "access_granted", "Opened dashboard page"),
    @("2025-12-25 00:35:26", "Admin", "quotation", "access_granted", "Opened quotation page"),
    @("2025-12-25 00:35:34", "Admin", "invoice", "access_granted", "Opened invoice page")
)

$startRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row + 1
for ($i = 0; $i -lt $logRows.Count; $i++) {
    $r = $startRow + $i
    $row = $logRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}
